$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column widths (values chosen so the saved OOXML column width
# is as close as possible to the target character widths)
$ws.Columns.Item(1).ColumnWidth = 11.3333333333333
$ws.Columns.Item(2).ColumnWidth = 20.6666666666667
$ws.Columns.Item(3).ColumnWidth = 21.5
$ws.Columns.Item(4).ColumnWidth = 19.3333333333333
$ws.Columns.Item(5).ColumnWidth = 17.3333333333333

# New / changed cell values in column E
$ws.Range("E3").Value = 144.342361450195
$ws.Range("E4").Value = 13400
$ws.Range("E5").Value = 134000
$ws.Range("E6").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("E9").Value = 0

# Update the selection to match the saved state
[void]$ws.Range("D21:E21").Select()

# Update the workbook window size/position
$win = $wb.Windows.Item(1)
$win.Left = -11025
$win.Top = 0
$win.Width = 11025
$win.Height = 10800
